$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (D1:I1), copying the header style (s="1") from C1
$headers = @("Processing Time", "ROC AUC", "Memory Usage", "Precision", "Accuracy", "Recall")
$cols = @("D", "E", "F", "G", "H", "I")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
}
# Copy style from C1 (existing header cell) to the new header cells D1:I1
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:I1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Row 2 - Naive Bayes
$ws.Range("B2").Value = "Naive Bayes"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.0206162929534912
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5703125
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1

# Row 3 - SVM
$ws.Range("B3").Value = "SVM"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.009204864501953101
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1

# Row 4 - MLP (name unchanged)
$ws.Range("C4").Value = 0.9665831244778612
$ws.Range("D4").Value = 0.1120295524597168
$ws.Range("E4").Value = 0.9933333333333332
$ws.Range("F4").Value = 0.88671875
$ws.Range("G4").Value = 0.9696969696969696
$ws.Range("H4").Value = 0.9666666666666668
$ws.Range("I4").Value = 0.9666666666666668

# Row 5 - DecisionTree
$ws.Range("B5").Value = "DecisionTree"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0.0281538963317871
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.7734375
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1

# Row 6 - RandomForest
$ws.Range("B6").Value = "RandomForest"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0.2764220237731933
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1.5625
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1

# Row 7 - KNN
$ws.Range("B7").Value = "KNN"
$ws.Range("C7").Value = 0.9665831244778612
$ws.Range("D7").Value = 0.0200996398925781
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.78125
$ws.Range("G7").Value = 0.9696969696969696
$ws.Range("H7").Value = 0.9666666666666668
$ws.Range("I7").Value = 0.9666666666666668

# Row 8 - LogReg
$ws.Range("B8").Value = "LogReg"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0.0265321731567382
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1

# Row 9 - GradientBoost
$ws.Range("B9").Value = "GradientBoost"
$ws.Range("C9").Value = 0.9665831244778612
$ws.Range("D9").Value = 0.1955735683441162
$ws.Range("E9").Value = 0.9741666666666668
$ws.Range("F9").Value = 1.12890625
$ws.Range("G9").Value = 0.9696969696969696
$ws.Range("H9").Value = 0.9666666666666668
$ws.Range("I9").Value = 0.9666666666666668

# Row 10 - XGBoost
$ws.Range("B10").Value = "XGBoost"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0.077242374420166
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 4.50390625
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1

# Row 11 - Custom AI Model
$ws.Range("B11").Value = "Custom AI Model"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0.0213708877563476
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.7109375
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
